$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.948.39"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "2.441.93"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.22"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.51"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +0.94%  "

$ws.Range("D9").Value = "2.432.80"
$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.05"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").Value = "2.901.44"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").Value = "62.990.81"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "2.446.67"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.31"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  +6.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.61"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  +12.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.33"
$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "613.87"
$ws.Range("E26").Value = "  +6.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.86"
$ws.Range("E27").Value = "  +1.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000103"
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("D29").Value = "2.604.80"
$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("E30").Value = "  +3.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").Value = "  -2.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.90"
$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("E34").Value = "  -4.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.18"
$ws.Range("E35").Value = "  +6.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").Value = "  -2.37%  "

$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.378"
$ws.Range("E38").Value = "  -1.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.70"
$ws.Range("E40").Value = "  -0.62%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "145.45"
$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("E42").Value = "  +10.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.78"
$ws.Range("E43").Value = "  -2.62%  "

$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.47"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.74"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.15"
$ws.Range("E47").Value = "  +3.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0535"
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.598"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0234"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0912"
$ws.Range("E51").Value = "  -1.29%  "
